# Actualización automática 2025-11-03 08:30:05
#
# This script reproduces the monthly "roll-forward" update applied to the
# commission workbook:
#   - "VENTAS POR GRUPO": the current-period (most recent month) figures are
#     zeroed out (a new, still-empty period has started) and the "x de 16"
#     progress labels in the totals row are reset to "0 de 16".
#   - "VENTA MENSUAL": the month headers shift forward by one month
#     (julio/agosto/septiembre/octubre -> agosto/septiembre/octubre/noviembre),
#     the per-client monthly figures shift along with them, and the column
#     widths are resized to match the new month columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": zero out the current month's figures
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$zeroCells = @("H3", "H8", "I8", "M8", "M9", "H10", "M10", "M11", "K12", "M12", "H15", "I15", "M15", "E16", "M16")
foreach ($cellRef in $zeroCells) {
    $wsGrupo.Range($cellRef).Value = 0
}

$resetLabels = @("E18", "H18", "I18", "K18", "M18")
foreach ($cellRef in $resetLabels) {
    $wsGrupo.Range($cellRef).Value = "0 de 16"
}

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": shift months forward by one column (C<-D<-E<-F),
# the last column becomes the new (empty) month.
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Month header labels (row 1, columns C:F). G1 "PRESUPUESTO" stays put.
$wsMensual.Range("C1").Value = "agosto"
$wsMensual.Range("D1").Value = "septiembre"
$wsMensual.Range("E1").Value = "octubre"
$wsMensual.Range("F1").Value = "noviembre"

# Column widths follow the same shift (new col G width stays 17).
$wsMensual.Columns.Item(3).ColumnWidth = 11.166666666666666  # -> 12
$wsMensual.Columns.Item(4).ColumnWidth = 15.166666666666666  # -> 16
$wsMensual.Columns.Item(5).ColumnWidth = 12.166666666666666  # -> 13
$wsMensual.Columns.Item(6).ColumnWidth = 14.166666666666666  # -> 15

# Per-client monthly data rows: each figure moves one column to the left
# (D takes E's old value, E takes F's old value), leaving column C as-is
# (already all zero) and setting F to 0 for the new empty month.
$dataRows = @(3, 8, 9, 10, 11, 12, 13, 15, 16, 18)
foreach ($row in $dataRows) {
    $eVal = $wsMensual.Cells.Item($row, 5).Value2
    $fVal = $wsMensual.Cells.Item($row, 6).Value2

    $wsMensual.Cells.Item($row, 4).Value = $eVal
    $wsMensual.Cells.Item($row, 5).Value = $fVal
    $wsMensual.Cells.Item($row, 6).Value = 0
}
